## Geral2014_17_d.xlsx - "Melhorias no projeto de faturamento, novo grafico"
##
## 1. Append the 201707 (July/2017) accumulation rows (731-742) to the
##    "Exportar Planilha" export sheet, mirroring the existing ANOMES_FT /
##    ANOMES_AT / QP / VA / VR layout.
## 2. Extend the source SQL (sheet "SQL") upper date bound from 201706 to
##    201707 so it matches the newly exported data.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "Exportar Planilha" - append rows 731..742
# ---------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Exportar Planilha")

$lastRow = 730

# New rows: ANOMES_FT, ANOMES_AT, QP, VA, VR
$newRows = @(
    @("201707","201606",9.0,217.4,9272.0),
    @("201707","201609",28.0,3588.49,12859.25),
    @("201707","201610",6.0,2103.35,-17.6),
    @("201707","201611",49.0,20505.7,24734.9),
    @("201707","201612",23.0,29310.71,23727.26),
    @("201707","201701",340.0,5514.53,27433.4),
    @("201707","201702",31.0,1195.99,57195.55),
    @("201707","201703",9332.0,189344.93,228220.14),
    @("201707","201704",39112.0,1336889.56,1105104.39),
    @("201707","201705",536619.0,7027827.35,855620.66),
    @("201707","201706",3553504.0,49039619.33,170144.74),
    @("201707","201707",2828194.0,69248014.87,0.0)
)

$firstNewRow = $lastRow + 1
$lastNewRow  = $lastRow + $newRows.Count

# Seed the new rows with the same formatting used by the row above them
# (column A/B plain text, C/D/E right-aligned "Dialog" numeric style) so
# the appended block matches the rest of the table.
$ws.Range("A$($lastRow):E$($lastRow)").Copy()
$ws.Range("A$($firstNewRow):E$($lastNewRow)").PasteSpecial(-4122)  # xlPasteFormats

# Force text storage for the ANOMES_FT / ANOMES_AT columns while we write
# the (numeric-looking) period codes, so they land as strings rather than
# numbers.
$ws.Range("A$($firstNewRow):B$($lastNewRow)").NumberFormat = "@"

$r = $firstNewRow
foreach ($row in $newRows) {
    $ws.Cells.Item($r, 1).Value2 = $row[0]
    $ws.Cells.Item($r, 2).Value2 = $row[1]
    $ws.Cells.Item($r, 3).Value2 = $row[2]
    $ws.Cells.Item($r, 4).Value2 = $row[3]
    $ws.Cells.Item($r, 5).Value2 = $row[4]
    $r = $r + 1
}

# Re-apply the reference formatting so the temporary "@" text format above
# doesn't linger on the new cells (the values are already strings, this
# just restores the normal/default display format to match the rest of
# the sheet).
$ws.Range("A$($lastRow):E$($lastRow)").Copy()
$ws.Range("A$($firstNewRow):E$($lastNewRow)").PasteSpecial(-4122)  # xlPasteFormats

# ---------------------------------------------------------------------
# 2) "SQL" - bump the upper bound of the exported period from 201706 to
#    201707
# ---------------------------------------------------------------------
$wsSql = $wb.Worksheets.Item("SQL")
$sqlCell = $wsSql.Cells.Item(2, 1)
$sql = $sqlCell.Value2
$sql = $sql.Replace("between '201401' and '201706'", "between '201401' and '201707'")
$sqlCell.Value2 = $sql

Write-Host "Added rows $($firstNewRow):$($lastNewRow) to 'Exportar Planilha' and updated SQL date range."
